# Applies the cryptos.xlsx cell content updates described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.886.72'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.622.66'
$ws.Range("E3").Value = '  -0.44%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.28'
$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("E6").Value = '  -1.52%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.96'
$ws.Range("E8").Value = '  -1.63%  '

$ws.Range("E9").Value = '  +0.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0605'
$ws.Range("E10").Value = '  -1.42%  '

$ws.Range("E11").Value = '  +0.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.58'
$ws.Range("E12").Value = '  -0.35%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.636.42'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("E14").Value = '  -0.93%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.551'
$ws.Range("E15").Value = '  -1.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.48'
$ws.Range("E16").Value = '  -1.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.891.05'
$ws.Range("E17").Value = '  -0.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '227.07'
$ws.Range("E18").Value = '  -1.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.58'
$ws.Range("E19").Value = '  -0.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0714'
$ws.Range("E20").Value = '  -1.12%  '

$ws.Range("E21").Value = '  +0.24%  '

$ws.Range("E22").Value = '  -0.40%  '

$ws.Range("E24").Value = '  +1.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.38'
$ws.Range("E25").Value = '  -0.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.89'
$ws.Range("E26").Value = '  -0.56%  '

$ws.Range("E27").Value = '  -0.66%  '

$ws.Range("E28").Value = '  +0.23%  '

$ws.Range("E29").Value = '  -1.30%  '

$ws.Range("E30").Value = '  -0.53%  '

$ws.Range("E31").Value = '  -0.49%  '

$ws.Range("E32").Value = '  -0.58%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.413.03'
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.08'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.61'
$ws.Range("E35").Value = '  +1.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.993'
$ws.Range("E36").Value = '  -2.40%  '

$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("E38").Value = '  -0.77%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.553'
$ws.Range("E39").Value = '  -0.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.846'
$ws.Range("E40").Value = '  -1.90%  '

$ws.Range("E41").Value = '  +0.11%  '

$ws.Range("E42").Value = '  -1.93%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.24'
$ws.Range("E43").Value = '  -1.64%  '

$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.81'
$ws.Range("E44").Value = '  -2.05%  '

$ws.Range("E45").Value = '  -2.00%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.764.79'
$ws.Range("E46").Value = '  -0.37%  '

$ws.Range("E47").Value = '  -3.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '89.04'
$ws.Range("E48").Value = '  +1.22%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.100'
$ws.Range("E49").Value = '  -0.32%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0502'
$ws.Range("E50").Value = '  -0.50%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.51'
$ws.Range("E51").Value = '  -0.63%  '
